# "filenames and class names correction for the sorting algos"
#
# The underlying data-generation run was repeated (new random timings) and a
# header row + a new "heap" column were added to the results sheet. This
# shifts the original 4 data rows down by one (now rows 2-5) and pushes the
# used range out to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new header row above the existing data (shifts rows down).
# ---------------------------------------------------------------------
$ws.Rows("1:1").Insert()

# ---------------------------------------------------------------------
# 2. Populate the header row, including the new "heap" column (G).
# ---------------------------------------------------------------------
$ws.Range("A1").Value2 = "n"
$ws.Range("B1").Value2 = "selection"
$ws.Range("C1").Value2 = "insertion"
$ws.Range("D1").Value2 = "merge"
$ws.Range("E1").Value2 = "quick"
$ws.Range("F1").Value2 = "count"
$ws.Range("G1").Value2 = "heap"

# ---------------------------------------------------------------------
# 3. Refresh the timing data (column A is unchanged; B:F rewritten with the
#    results of the new benchmark run).
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = 10
$ws.Range("A3").Value2 = 100
$ws.Range("A4").Value2 = 1000
$ws.Range("A5").Value2 = 10000

$ws.Range("B2").Value2 = [double]"4.2611476167754799E-5"
$ws.Range("B3").Value2 = [double]"8.1452718668071104E-4"
$ws.Range("B4").Value2 = [double]"7.5940370990121603E-2"
$ws.Range("B5").Value2 = [double]"7.8670451750455497"

$ws.Range("C2").Value2 = [double]"1.5758811410097401E-5"
$ws.Range("C3").Value2 = [double]"9.2170824232867604E-4"
$ws.Range("C4").Value2 = [double]"9.4620218767875505E-2"
$ws.Range("C5").Value2 = [double]"9.4419660158097791"

$ws.Range("D2").Value2 = [double]"4.0384687604509001E-5"
$ws.Range("D3").Value2 = [double]"6.0662312667096E-4"
$ws.Range("D4").Value2 = [double]"7.7411589423517398E-3"
$ws.Range("D5").Value2 = [double]"9.9334049540454702E-2"

$ws.Range("E2").Value2 = [double]"2.35215931540189E-5"
$ws.Range("E3").Value2 = [double]"3.34251533138165E-4"
$ws.Range("E4").Value2 = [double]"5.2356828647174003E-3"
$ws.Range("E5").Value2 = [double]"0.132511176510582"

$ws.Range("F2").Value2 = [double]"4.5158980603550099E-5"
$ws.Range("F3").Value2 = [double]"9.2227678405723697E-5"
$ws.Range("F4").Value2 = [double]"4.5999766401982401E-4"
$ws.Range("F5").Value2 = [double]"4.0754531307029499E-3"

# ---------------------------------------------------------------------
# 4. Update the worksheet-scoped defined name "results" to the new extent.
# ---------------------------------------------------------------------
$nm = $wb.Names.Item(1)
$nm.RefersTo = "=Sheet1!`$A`$2:`$F`$5"

# ---------------------------------------------------------------------
# 5. Re-point the chart series at the shifted ranges.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$chart.SeriesCollection().Item(1).Formula = '=SERIES("selection",Sheet1!$A$2:$A$5,Sheet1!$B$2:$B$5,1)'
$chart.SeriesCollection().Item(2).Formula = '=SERIES("insertion",Sheet1!$A$2:$A$5,Sheet1!$C$2:$C$5,2)'
$chart.SeriesCollection().Item(3).Formula = '=SERIES("merging",Sheet1!$A$2:$A$5,Sheet1!$D$2:$D$5,3)'
$chart.SeriesCollection().Item(4).Formula = '=SERIES("quick",Sheet1!$A$2:$A$5,Sheet1!$E$2:$E$5,4)'
$chart.SeriesCollection().Item(5).Formula = '=SERIES("count",Sheet1!$A$2:$A$5,Sheet1!$F$2:$F$5,5)'

# ---------------------------------------------------------------------
# 6. Move/resize the chart to its new anchor position.
# ---------------------------------------------------------------------
$co.Left = 462.6874803149606
$co.Top = 1.5
$co.Width = 756.6874015748033
$co.Height = 312.0

# ---------------------------------------------------------------------
# 7. Restore the active selection to the cell the author left selected.
# ---------------------------------------------------------------------
$ws.Range("G3").Select()
